$d = $word.ActiveDocument

function Get-Text {
    return $d.Content.Text
}

# ---------------------------------------------------------------------------
# Change 1: merge "()" run with ": Returns a collection of " run into a single
# run reading "(): Returns a collection of ", while leaving the following
# "records from the history file." run untouched (separate run).
# ---------------------------------------------------------------------------

# Protect the boundary so the merge doesn't sweep up the next run too: toggle
# a formatting-only change (no text change) on "records from the history
# file." - formatting-only edits never trigger a run merge in this engine.
$text = Get-Text
$guard1 = "records from the history file."
$gIdx = $text.IndexOf($guard1)
$gRng = $d.Range($gIdx, $gIdx + $guard1.Length)
$gRng.Bold = 1

# Now merge "()" and ": Returns a collection of " together. Add a throwaway
# trailing marker character so the replacement text differs from the
# original (identical-text assignments are no-ops in this engine and would
# not actually merge the runs).
$text = Get-Text
$mergeTarget = "()" + ": Returns a collection of "
$mIdx = $text.IndexOf($mergeTarget)
$mRng = $d.Range($mIdx, $mIdx + $mergeTarget.Length)
$mRng.Text = "(): Returns a collection of ~"

# Remove the throwaway marker with a small isolated edit (this does not
# reach into the guarded run, so no further unwanted merging happens).
$text = Get-Text
$tmpMarker = "(): Returns a collection of ~"
$tIdx = $text.IndexOf($tmpMarker)
$tRng = $d.Range($tIdx, $tIdx + $tmpMarker.Length)
$tRng.Text = "(): Returns a collection of "

# Remove the formatting guard now that the merge boundary is settled.
$text = Get-Text
$gIdx2 = $text.IndexOf($guard1)
$gRng2 = $d.Range($gIdx2, $gIdx2 + $guard1.Length)
$gRng2.Bold = 0

# ---------------------------------------------------------------------------
# Change 2: merge " e)" with ": Commit a sell request." into a single run
# reading " e): Commit a sell request." (end of paragraph, nothing after it
# to worry about).
# ---------------------------------------------------------------------------

$text = Get-Text
$sellMarker = "SellButton_Click(object sender, RoutedEventArgs" + " e): Commit a sell request."
$sIdx = $text.IndexOf($sellMarker)
$sStart = $sIdx + ("SellButton_Click(object sender, RoutedEventArgs").Length
$sSub = " e): Commit a sell request."
$sRng = $d.Range($sStart, $sStart + $sSub.Length)
$sRng.Text = " e): Commit a sell request.~"

$text = Get-Text
$sellTmp = " e): Commit a sell request.~"
$sIdx2 = $text.IndexOf($sellTmp)
$sRng2 = $d.Range($sIdx2, $sIdx2 + $sellTmp.Length)
$sRng2.Text = " e): Commit a sell request."

# ---------------------------------------------------------------------------
# Change 3: rename the 3rd Button_Click handler (the cancel/delete one) to
# CancelRequestButton_Click.
# ---------------------------------------------------------------------------

$text = Get-Text
$renameMarker = "Button_Click(object sender, RoutedEventArgs e): Deletes a pending"
$rIdx = $text.IndexOf($renameMarker)
$rRng = $d.Range($rIdx, $rIdx + ("Button_Click").Length)
$rRng.Text = "CancelRequestButton_Click"

# ---------------------------------------------------------------------------
# Move the hidden "_GoBack" bookmark from after "Deletes a pending from the
# server." to right after the " e)" of the CancelRequestButton_Click
# signature (before the ": Deletes a pending..." text).
# ---------------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$text = Get-Text
$goBackMarker = "CancelRequestButton_Click(object sender, RoutedEventArgs e)"
$gbIdx = $text.IndexOf($goBackMarker)
$gbPos = $gbIdx + $goBackMarker.Length
$gbRng = $d.Range($gbPos, $gbPos)
$d.Bookmarks.Add("_GoBack", $gbRng)

